$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-sort the pair of same-kickoff-time matches at rows 114/115: ---
# --- swap every column except A (rank) to match the refreshed source order ---
$ws.Range("B114").Value = 7559468
$ws.Range("F114").Value = "Liverpool Montevideo"
$ws.Range("G114").Value = "CA River Plate"
$ws.Range("H114").Value = 2
$ws.Range("I114").Value = 1
$ws.Range("J114").Value = "H"
$ws.Range("K114").Value = 1.7
$ws.Range("L114").Value = 3
$ws.Range("M114").Value = 5.75
$ws.Range("N114").Value = 1.833
$ws.Range("P114").Value = 4.5
$ws.Range("Q114").Value = -0.5
$ws.Range("R114").Value = 1.925
$ws.Range("S114").Value = 1.925
$ws.Range("T114").Value = 2.25
$ws.Range("U114").Value = 2.025
$ws.Range("V114").Value = 1.825
$ws.Range("W114").Value = 0.833
$ws.Range("X114").Value = -1
$ws.Range("Z114").Value = 0.925
$ws.Range("AA114").Value = -1
$ws.Range("AB114").Value = 1.025
$ws.Range("AC114").Value = -1

$ws.Range("B115").Value = 7559469
$ws.Range("F115").Value = "Montevideo Wanderers"
$ws.Range("G115").Value = "Penarol"
$ws.Range("H115").Value = 0
$ws.Range("I115").Value = 0
$ws.Range("J115").Value = "D"
$ws.Range("K115").Value = 4.75
$ws.Range("L115").Value = 3.4
$ws.Range("M115").Value = 1.7
$ws.Range("N115").Value = 2.7
$ws.Range("P115").Value = 2.45
$ws.Range("Q115").Value = 0
$ws.Range("R115").Value = 2.05
$ws.Range("S115").Value = 1.8
$ws.Range("T115").Value = 2.5
$ws.Range("U115").Value = 1.975
$ws.Range("V115").Value = 1.875
$ws.Range("W115").Value = -1
$ws.Range("X115").Value = 2.2
$ws.Range("Z115").Value = 0
$ws.Range("AA115").Value = 0
$ws.Range("AB115").Value = -1
$ws.Range("AC115").Value = 0.875

# --- Same re-sort for the pair at rows 118/119 ---
$ws.Range("B118").Value = 7013885
$ws.Range("F118").Value = "La Luz"
$ws.Range("G118").Value = "Atletico Fenix Montevideo"
$ws.Range("K118").Value = 3
$ws.Range("L118").Value = 3
$ws.Range("M118").Value = 2.4
$ws.Range("N118").Value = 2.9
$ws.Range("O118").Value = 2.75
$ws.Range("P118").Value = 2.6
$ws.Range("Q118").Value = 0
$ws.Range("R118").Value = 2.025
$ws.Range("S118").Value = 1.825
$ws.Range("T118").Value = 2
$ws.Range("U118").Value = 2.025
$ws.Range("V118").Value = 1.825
$ws.Range("Y118").Value = 1.6
$ws.Range("AA118").Value = 0.825
$ws.Range("AB118").Value = 0
$ws.Range("AC118").Value = 0

$ws.Range("B119").Value = 7013702
$ws.Range("F119").Value = "Defensor Sporting"
$ws.Range("G119").Value = "Danubio"
$ws.Range("K119").Value = 1.8
$ws.Range("L119").Value = 3.6
$ws.Range("M119").Value = 4.2
$ws.Range("N119").Value = 1.8
$ws.Range("O119").Value = 3.6
$ws.Range("P119").Value = 4.2
$ws.Range("Q119").Value = -0.75
$ws.Range("R119").Value = 2.05
$ws.Range("S119").Value = 1.8
$ws.Range("T119").Value = 2.25
$ws.Range("U119").Value = 1.85
$ws.Range("V119").Value = 2
$ws.Range("Y119").Value = 3.2
$ws.Range("AA119").Value = 0.8
$ws.Range("AB119").Value = -0.5
$ws.Range("AC119").Value = 0.5

# --- Stamp cell formatting (style + number format) for the 3 brand-new rows ---
# --- by cloning row 165 (which already carries the correct A/E column styles) ---
$ws.Range("A165:AC165").Copy($ws.Range("A166:AC166"))
$ws.Range("A165:AC165").Copy($ws.Range("A167:AC167"))
$ws.Range("A165:AC165").Copy($ws.Range("A168:AC168"))

# --- Row 165 now becomes a newly-finished match (odds/result refreshed) ---
$ws.Range("B165").Value = 7994684
$ws.Range("E165").Value = 45381.79166666666
$ws.Range("F165").Value = "Defensor Sporting"
$ws.Range("G165").Value = "Danubio"
$ws.Range("H165").Value = 1
$ws.Range("I165").Value = 0
$ws.Range("J165").Value = "H"
$ws.Range("K165").Value = 1.909
$ws.Range("L165").Value = 3.25
$ws.Range("M165").Value = 3.8
$ws.Range("N165").Value = 1.75
$ws.Range("O165").Value = 3.3
$ws.Range("P165").Value = 4.5
$ws.Range("Q165").Value = -0.75
$ws.Range("R165").Value = 2.025
$ws.Range("S165").Value = 1.825
$ws.Range("U165").Value = 1.975
$ws.Range("V165").Value = 1.875
$ws.Range("W165").Value = 0.75
$ws.Range("X165").Value = -1
$ws.Range("Y165").Value = -1
$ws.Range("Z165").Value = 0.5125
$ws.Range("AA165").Value = -0.5
$ws.Range("AB165").Value = -1
$ws.Range("AC165").Value = 0.875

# --- New row 166: Club Atletico Progreso vs Deportivo Maldonado ---
$ws.Range("A166").Value = 164
$ws.Range("B166").Value = 7995146
$ws.Range("C166").Value = "Uruguay Primera División"
$ws.Range("D166").Value = "Uruguay Apertura"
$ws.Range("E166").Value = 45382.625
$ws.Range("F166").Value = "Club Atletico Progreso"
$ws.Range("G166").Value = "Deportivo Maldonado"
$ws.Range("H166").Value = 3
$ws.Range("I166").Value = 1
$ws.Range("J166").Value = "H"
$ws.Range("K166").Value = 2.4
$ws.Range("L166").Value = 3.1
$ws.Range("M166").Value = 3
$ws.Range("N166").Value = 2.15
$ws.Range("O166").Value = 3.25
$ws.Range("P166").Value = 3.1
$ws.Range("Q166").Value = -0.25
$ws.Range("R166").Value = 1.925
$ws.Range("S166").Value = 1.925
$ws.Range("T166").Value = 2.5
$ws.Range("U166").Value = 1.975
$ws.Range("V166").Value = 1.875
$ws.Range("W166").Value = 1.15
$ws.Range("X166").Value = -1
$ws.Range("Y166").Value = -1
$ws.Range("Z166").Value = 0.925
$ws.Range("AA166").Value = -1
$ws.Range("AB166").Value = 0.9750000000000001
$ws.Range("AC166").Value = -1

# --- New row 167: Miramar Misiones vs Cerro ---
$ws.Range("A167").Value = 165
$ws.Range("B167").Value = 7995141
$ws.Range("C167").Value = "Uruguay Primera División"
$ws.Range("D167").Value = "Uruguay Apertura"
$ws.Range("E167").Value = 45382.72916666666
$ws.Range("F167").Value = "Miramar Misiones"
$ws.Range("G167").Value = "Cerro"
$ws.Range("H167").Value = 2
$ws.Range("I167").Value = 2
$ws.Range("J167").Value = "D"
$ws.Range("K167").Value = 2.6
$ws.Range("L167").Value = 3
$ws.Range("M167").Value = 2.75
$ws.Range("N167").Value = 2.5
$ws.Range("O167").Value = 3
$ws.Range("P167").Value = 2.875
$ws.Range("Q167").Value = 0
$ws.Range("R167").Value = 1.8
$ws.Range("S167").Value = 2.05
$ws.Range("T167").Value = 2.25
$ws.Range("U167").Value = 2
$ws.Range("V167").Value = 1.85
$ws.Range("W167").Value = -1
$ws.Range("X167").Value = 2
$ws.Range("Y167").Value = -1
$ws.Range("Z167").Value = 0
$ws.Range("AA167").Value = 0
$ws.Range("AB167").Value = 1
$ws.Range("AC167").Value = -1

# --- New row 168: Montevideo Wanderers vs Boston River (was pending row 165) ---
$ws.Range("A168").Value = 166
$ws.Range("B168").Value = 7994683
$ws.Range("C168").Value = "Uruguay Primera División"
$ws.Range("D168").Value = "Uruguay Apertura"
$ws.Range("E168").Value = 45382.83333333334
$ws.Range("F168").Value = "Montevideo Wanderers"
$ws.Range("G168").Value = "Boston River"
$ws.Range("H168").Value = 0
$ws.Range("I168").Value = 2
$ws.Range("J168").Value = "A"
$ws.Range("K168").Value = 2.5
$ws.Range("L168").Value = 3.1
$ws.Range("M168").Value = 2.75
$ws.Range("N168").Value = 2.8
$ws.Range("O168").Value = 3.1
$ws.Range("P168").Value = 2.5
$ws.Range("Q168").Value = 0
$ws.Range("R168").Value = 2.05
$ws.Range("S168").Value = 1.8
$ws.Range("T168").Value = 2.25
$ws.Range("U168").Value = 2
$ws.Range("V168").Value = 1.85
$ws.Range("W168").Value = -1
$ws.Range("X168").Value = -1
$ws.Range("Y168").Value = 1.5
$ws.Range("Z168").Value = -1
$ws.Range("AA168").Value = 0.8
$ws.Range("AB168").Value = -0.5
$ws.Range("AC168").Value = 0.425

